# BOT; UPDATE DATA (#1655)
# Applies the data refresh + view-state changes described by the diff:
#  - "all" sheet becomes the active/selected tab (was "other")
#  - updated daily totals (column H) on "all" sheet, rows 26-41 (+ C40)
#  - updated daily totals (column J, plus E82) on "kobe" sheet, rows 81-96
#  - refreshed cell selections on "all", "kobe" and "other" sheets

$wb = $excel.ActiveWorkbook

$wsAll   = $wb.Worksheets.Item("all")
$wsKobe  = $wb.Worksheets.Item("kobe")
$wsOther = $wb.Worksheets.Item("other")

# ---------------------------------------------------------------------------
# "all" sheet - updated figures (column H, rows 26-41; C40)
# ---------------------------------------------------------------------------
$wsAll.Range("H26").Value = 165
$wsAll.Range("H27").Value = 165
$wsAll.Range("H28").Value = 176
$wsAll.Range("H29").Value = 179
$wsAll.Range("H30").Value = 184
$wsAll.Range("H31").Value = 189
$wsAll.Range("H32").Value = 196
$wsAll.Range("H33").Value = 198
$wsAll.Range("H34").Value = 202
$wsAll.Range("H35").Value = 213
$wsAll.Range("H36").Value = 213
$wsAll.Range("H37").Value = 219
$wsAll.Range("H38").Value = 223
$wsAll.Range("H39").Value = 227
$wsAll.Range("C40").Value = 281
$wsAll.Range("H40").Value = 228
$wsAll.Range("H41").Value = 229

# ---------------------------------------------------------------------------
# "kobe" sheet - updated figures (column J, rows 81-96; E82)
# ---------------------------------------------------------------------------
$wsKobe.Range("J81").Value = 158
$wsKobe.Range("E82").Value = 268
$wsKobe.Range("J82").Value = 158
$wsKobe.Range("J83").Value = 169
$wsKobe.Range("J84").Value = 172
$wsKobe.Range("J85").Value = 177
$wsKobe.Range("J86").Value = 180
$wsKobe.Range("J87").Value = 187
$wsKobe.Range("J88").Value = 189
$wsKobe.Range("J89").Value = 193
$wsKobe.Range("J90").Value = 204
$wsKobe.Range("J91").Value = 204
$wsKobe.Range("J92").Value = 210
$wsKobe.Range("J93").Value = 214
$wsKobe.Range("J94").Value = 218
$wsKobe.Range("J95").Value = 219
$wsKobe.Range("J96").Value = 220

# ---------------------------------------------------------------------------
# Refresh each sheet's remembered selection (bottom-right pane of the frozen
# view) to match the saved workbook state.
# ---------------------------------------------------------------------------
$wsKobe.Activate()
$wsKobe.Range("A80").Select()

$wsOther.Activate()
$wsOther.Range("A72").Select()

# "all" is the sheet that ends up active/selected when the workbook is saved.
$wsAll.Activate()
$wsAll.Range("L35").Select()
